# The "Fallbeispiel" question block (rows 10-15: the prompt row plus its
# five operationalized-goal sub-bullets) was removed from the question
# list. Deleting the entire rows shifts everything below up by six rows,
# which is exactly what the canonical XML diff shows (row N after the
# edit == row N+6 before, same per-row style, same shared-string text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10:A15").EntireRow.Delete()

# After the deletion Excel leaves the selection on the row that slid up
# into the gap (previously row 16, now row 10's neighbourhood settles at
# the former row 26 "Nennen Sie mindestens vier Inhalte..." -> A20).
$ws.Range("A20").Select()
